$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45188 -> 45189) for every data row (rows 2 through 271).
for ($r = 2; $r -le 271; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}
